$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 136
$ws.Range("D2").Value = 59657.42
$ws.Range("E2").Value = 117.78

# Row 3
$ws.Range("C3").Value = 55
$ws.Range("D3").Value = 50110
$ws.Range("E3").Value = 955

# Row 4
$ws.Range("C4").Value = 56
$ws.Range("D4").Value = 47130
$ws.Range("E4").Value = 985

# Row 5
$ws.Range("C5").Value = 68
$ws.Range("D5").Value = 46247.31
$ws.Range("E5").Value = 673.4

# Row 6
$ws.Range("C6").Value = 63
$ws.Range("D6").Value = 45480
$ws.Range("E6").Value = 680

# Row 7
$ws.Range("C7").Value = 55
$ws.Range("D7").Value = 43600
$ws.Range("E7").Value = 730

# Row 8
$ws.Range("C8").Value = 68
$ws.Range("D8").Value = 40305
$ws.Range("E8").Value = 610

# Row 9
$ws.Range("C9").Value = 68
$ws.Range("D9").Value = 39140

# Row 10
$ws.Range("C10").Value = 68
$ws.Range("D10").Value = 33770
$ws.Range("E10").Value = 500

# Row 11
$ws.Range("C11").Value = 68
$ws.Range("D11").Value = 27235

# Row 12
$ws.Range("C12").Value = 68
$ws.Range("D12").Value = 25059.55
$ws.Range("E12").Value = 379.87

# Row 13
$ws.Range("C13").Value = 68
$ws.Range("D13").Value = 22479.27
$ws.Range("E13").Value = 330.96

# Row 14
$ws.Range("C14").Value = 68
$ws.Range("D14").Value = 14660.17
$ws.Range("E14").Value = 219.26

# Row 15
$ws.Range("C15").Value = 68
$ws.Range("D15").Value = 10026.17
$ws.Range("E15").Value = 156.52

# Row 16
$ws.Range("C16").Value = 68
$ws.Range("D16").Value = 8973.870000000001
$ws.Range("E16").Value = 139.51

# Row 17
$ws.Range("C17").Value = 68
$ws.Range("D17").Value = 8233.610000000001
$ws.Range("E17").Value = 121.3

# Row 18
$ws.Range("C18").Value = 68
$ws.Range("D18").Value = 7682.6
$ws.Range("E18").Value = 117.26

# Row 19
$ws.Range("C19").Value = 68
$ws.Range("D19").Value = 7607.76
$ws.Range("E19").Value = 117.09

# Row 20
$ws.Range("C20").Value = 68
$ws.Range("D20").Value = 7445.98
$ws.Range("E20").Value = 115.57

# Row 21
$ws.Range("C21").Value = 68
$ws.Range("D21").Value = 7337.1
$ws.Range("E21").Value = 114.24

# Row 22
$ws.Range("C22").Value = 68
$ws.Range("D22").Value = 7210.69
$ws.Range("E22").Value = 112.27

# Row 23
$ws.Range("C23").Value = 68
$ws.Range("D23").Value = 6801.74
$ws.Range("E23").Value = 95.36

# Row 24
$ws.Range("C24").Value = 68
$ws.Range("D24").Value = 6750.14
$ws.Range("E24").Value = 99.84999999999999

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 58.51
$ws.Range("E26").Value = -2.72

# Row 27
$ws.Range("B27").Value = 14
$ws.Range("D27").Value = 58
$ws.Range("E27").Value = 4.16

# Row 31
$ws.Range("B31").Value = 9
$ws.Range("D31").Value = 36.41
$ws.Range("E31").Value = 2.14

# Row 32
$ws.Range("A32").Value = 'BANK OF AFRICA ML (BOAM)'
$ws.Range("B32").Value = 11
$ws.Range("C32").Value = 8
$ws.Range("D32").Value = 28.86
$ws.Range("E32").Value = 2

# Row 33
$ws.Range("A33").Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws.Range("B33").Value = 8
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 28.57
$ws.Range("E33").Value = 1.48

# Row 37
$ws.Range("A37").Value = 'CIE CI (CIEC)'
$ws.Range("B37").Value = 7
$ws.Range("C37").Value = 4
$ws.Range("D37").Value = 20.65
$ws.Range("E37").Value = -2.78

# Row 38
$ws.Range("A38").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws.Range("B38").Value = 5
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 15.74
$ws.Range("E38").Value = 4.1

# Row 39
$ws.Range("C39").Value = 11
$ws.Range("D39").Value = 14.42
$ws.Range("E39").Value = -3.85

# Row 40
$ws.Range("A40").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws.Range("B40").Value = 11
$ws.Range("C40").Value = 10
$ws.Range("D40").Value = 14.37
$ws.Range("E40").Value = 4.25

# Row 41
$ws.Range("A41").Value = 'UNIWAX CI (UNXC)'
$ws.Range("B41").Value = 12
$ws.Range("C41").Value = 9
$ws.Range("D41").Value = 14.12
$ws.Range("E41").Value = 2.53

# Row 42
$ws.Range("A42").Value = 'PALM CI (PALC)'
$ws.Range("B42").Value = 8
$ws.Range("C42").Value = 6
$ws.Range("D42").Value = 14.07
$ws.Range("E42").Value = -7.48

# Row 44
$ws.Range("A44").Value = 'SMB CI (SMBC)'
$ws.Range("B44").Value = 12
$ws.Range("C44").Value = 12
$ws.Range("D44").Value = 10.96
$ws.Range("E44").Value = 1.13

# Row 45
$ws.Range("A45").Value = 'SUCRIVOIRE (SCRC)'
$ws.Range("B45").Value = 8
$ws.Range("C45").Value = 10
$ws.Range("D45").Value = 10.74
$ws.Range("E45").Value = -1.01

# Row 54
$ws.Range("A54").Value = 'TOTAL'
$ws.Range("B54").Value = 0
$ws.Range("C54").Value = 67
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = 0

# Row 55
$ws.Range("A55").Value = 'SOGB CI (SOGC)'
$ws.Range("B55").Value = 8
$ws.Range("C55").Value = 6
$ws.Range("D55").Value = -1.95
$ws.Range("E55").Value = -3.67

# Row 56
$ws.Range("A56").Value = 'SOLIBRA CI (SLBC)'
$ws.Range("B56").Value = 13
$ws.Range("C56").Value = 13
$ws.Range("D56").Value = -2.3
$ws.Range("E56").Value = 2.11

# Row 57
$ws.Range("A57").Value = 'SONATEL SN (SNTS)'
$ws.Range("B57").Value = 2
$ws.Range("C57").Value = 3
$ws.Range("D57").Value = -3.43
$ws.Range("E57").Value = 0.8

# Row 58
$ws.Range("A58").Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$ws.Range("B58").Value = 11
$ws.Range("C58").Value = 15
$ws.Range("D58").Value = -4.04
$ws.Range("E58").Value = -2.29

# Row 66
$ws.Range("A66").Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$ws.Range("C66").Value = 9
$ws.Range("D66").Value = -27.25
$ws.Range("E66").Value = -0.59

# Row 67
$ws.Range("A67").Value = 'SICOR CI (SICC)'
$ws.Range("C67").Value = 6
$ws.Range("D67").Value = -30.3
$ws.Range("E67").Value = -6.82
